$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new display text, as scraped from the OOXML diff.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (preserving formatting like trailing zeros / percent signs)
# instead of re-interpreting them as numbers, then the style is reset back
# to "Normal" so no stray number-format / quote-prefix style is introduced.
$updates = @{
    "D2" = '304.02'
    "E2" = '6.01%'
    "D3" = '32.44'
    "E3" = '11.43%'
    "D4" = '5.300'
    "E4" = '2.02%'
    "D5" = '0.07497'
    "E5" = '7.54%'
    "D6" = '7.848'
    "D7" = '3.805'
    "E7" = '6.93%'
    "D8" = '1.536'
    "E8" = '9.95%'
    "D9" = '0.9182'
    "E9" = '2.47%'
    "D10" = '0.01771'
    "E10" = '2,631.67%'
    "D11" = '0.1698'
    "E11" = '5.79%'
    "D12" = '0.07857'
    "E12" = '3.22%'
    "D13" = '0.07952'
    "E13" = '3.90%'
    "D14" = '0.03018'
    "E14" = '3.32%'
    "D15" = '0.09917'
    "E15" = '10.12%'
    "D16" = '0.001496'
    "E16" = '-5.96%'
    "D17" = '0.04615'
    "E17" = '1.96%'
    "D18" = '0.006275'
    "E18" = '-0.34%'
    "D19" = '3.476'
    "E19" = '0.53%'
    "D20" = '2.230'
    "E20" = '0.00%'
    "E21" = '3.06%'
    "D22" = '0.1342'
    "E22" = '0.80%'
    "D23" = '4.502'
    "E23" = '12.36%'
    "E24" = '1.40%'
    "D25" = '0.001219'
    "E25" = '0.93%'
    "D26" = '0.004458'
    "E26" = '5.28%'
    "E27" = '19.70%'
    "D28" = '0.0001933'
    "E28" = '19.24%'
    "D40" = '0.04536'
    "E40" = '5.13%'
    "D41" = '0.007182'
    "E41" = '3.55%'
    "D42" = '0.1349'
    "E42" = '8.39%'
    "D43" = '0.002206'
    "E43" = '6.29%'
    "D44" = '0.01282'
    "E44" = '11.52%'
    "D45" = '0.00006200'
    "E45" = '6.05%'
    "D46" = '0.7082'
    "E46" = '-63.29%'
    "D47" = '0.01297'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
